# Update the workbook to add the "phi4" results sheet and refresh a few
# values/labels on the existing "qwen2.5" sheet, per the submodule bump
# for LongGenBench_Test.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "qwen2.5"

# ---------------------------------------------------------------------
# qwen2.5 sheet: relabel row 2, refresh row 4's value, and append two
# new rows (h2o + the original snapkv row, now moved to the bottom).
# ---------------------------------------------------------------------
$ws1.Range("A2").Value = "preds_ns5_ws3_st20.0_ea1.0_snks0_hopf_False_type_max_fused_lenNone_gblFalse"
$ws1.Range("B2").Value = 9591.93088

$ws1.Range("B4").Value = 12283.930624

$ws1.Range("A2").Copy()
$ws1.Range("A5:A6").PasteSpecial(-4122)   # xlPasteFormats - reuse existing label style

$ws1.Range("A5").Value = "preds_ns5_ws200_st10.0_ea1.0_snks0_hopf_True_type_h2o_lenNone_gblFalse"
$ws1.Range("B5").Value = 17808.02688

$ws1.Range("A6").Value = "preds_ns5_ws200_st20.0_ea1.0_snks0_hopf_True_type_snapkv_lenNone_gblFalse"
$ws1.Range("B6").Value = 0

# ---------------------------------------------------------------------
# Add the new "phi4" sheet right after "qwen2.5" and populate it.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "phi4"

$ws1.Range("B1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws2.Range("B1").Value = "phi4"

$ws1.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

$ws2.Range("A2").Value = "preds_ns5_ws200_st20.0_ea1.0_snks0_hopf_True_type_sum_fused_lenNone_gblFalse"
$ws2.Range("B2").Value = 44093.696

$ws2.Range("A3").Value = "preds_ns5_ws200_st10.0_ea1.0_snks0_hopf_True_type_h2o_lenNone_gblFalse"
$ws2.Range("B3").Value = 36342.912

$ws2.Range("A4").Value = "preds_ns5_ws200_st20.0_ea1.0_snks0_hopf_True_type_max_fused_lenNone_gblFalse"
$ws2.Range("B4").Value = 44093.696

$ws2.Range("A5").Value = "preds_ns5_ws32_st1025.0_ea1.0_snks0_hopf_True_type_snapkv_lenNone_gblFalse"
$ws2.Range("B5").Value = 140582.0928

$excel.CutCopyMode = $false
$ws1.Select()
